$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Activate()

# New header cell for column K (copy the header style from J1, same row)
$ws.Range("K1").Value = "Order result"
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2 already has A2:J2 populated - just add the new K2 result cell
$ws.Range("K2").Value = "Mua;AAA;Chờ gửi;LO;100;25,000;"

# Row 3 - new order: Mua AAA qty 99
$ws.Range("A3").Value = "046FIA0016"
$ws.Range("B3").Value = 123
$ws.Range("C3").Value = "'0001000021"
$ws.Range("D3").Value = "Lệnh thông thường"
$ws.Range("E3").Value = "AAA"
$ws.Range("F3").Value = 99
$ws.Range("G3").Value = 25
$ws.Range("H3").Value = "Mua"
$ws.Range("I3").Value = "LO"
$ws.Range("J3").Value = 123
$ws.Range("K3").Value = "Mua;AAA;Chờ gửi;LO;99;25,000;"

# Row 4 (temporary) - new order: Bán AAA qty 100
$ws.Range("A4").Value = "046FIA0016"
$ws.Range("B4").Value = 123
$ws.Range("C4").Value = "'0001000021"
$ws.Range("D4").Value = "Lệnh thông thường"
$ws.Range("E4").Value = "AAA"
$ws.Range("F4").Value = 100
$ws.Range("G4").Value = 25
$ws.Range("H4").Value = "Bán"
$ws.Range("I4").Value = "LO"
$ws.Range("J4").Value = 123
$ws.Range("K4").Value = "Bán;AAA;Chờ gửi;LO;100;25,000;"

# Insert a new row above the row just entered - this pushes the
# Bán/AAA/100 order down to row 5, matching the final layout, and
# reproduces the shared-string insertion order seen in the real edit.
$ws.Rows.Item(4).Insert()

# Row 4 (final) - new order: Bán BSL qty 15
$ws.Range("A4").Value = "046FIA0016"
$ws.Range("B4").Value = 123
$ws.Range("C4").Value = "'0001000021"
$ws.Range("D4").Value = "Lệnh thông thường"
$ws.Range("E4").Value = "BSL"
$ws.Range("F4").Value = 15
$ws.Range("G4").Value = 11
$ws.Range("H4").Value = "Bán"
$ws.Range("I4").Value = "LO"
$ws.Range("J4").Value = 123
$ws.Range("K4").Value = "Bán;BSL;Chờ gửi;LO;15;11,000;"

# Widen the new result column
$ws.Columns.Item(11).ColumnWidth = 40.8

# Final selection cell
$ws.Range("E7").Select()
